$d = $word.ActiveDocument
$tables = $d.Tables

$texts = @(
    "במחצית זאת עשינו מלאכות לסוכה, התאמנו על סריגה ועשינו מלאכות יד ותכשיטים לתפארה.`nהייתה אוירה טובה בכיתה וכולן נהנו.`nגילה את ילדה נהדרת, בהצלחה!",
    "במחצית זאת למדנו את תורת המיספרים, הכרנו כל מספר לעומק, והיתחלנו עם פעולות חשבון בסיסיות,חיבור וחיסור, התקדמנו הרבה עם הספר ""חושבים 1"".והתכוננו לקראת השנה החדשה בההכרה מלמעלה כל כפל וחילוק,`nגילה, את ילדה נפלאה,עלי והצלחי!!",
    "במחצית זאת למדנו על עולם המוזיקה, על התווים ועל רמות הקול, התעסקנו עם שירים על מעגל השנה, הייתה אוירה כיפית ונחמדה.`nגילה את תלמידה מדהימה!",
    "במחצית זאת למדנו חומש בראשית, למדנו והתפעלנו מבריאת העולם, עקידת יצחק וכו....`nגילה את תלמידה מצויוינת, בהצלחה!"
)

for ($i = 1; $i -le $tables.Count; $i++) {
    $cell = $tables.Item($i).Cell(1, 2)
    $cell.Range.Text = $texts[$i - 1]
}
